$d = $word.ActiveDocument

# The document ends with the "HandmadeList" lint block:
#   @BEGIN HandmadeList
#   Списки должны быть оформлены средствами Word, не вручную.
#   @END
#
# We need to insert a new lint block right before that closing "@END",
# i.e. end up with:
#   @BEGIN HandmadeList
#   Списки должны быть оформлены средствами Word, не вручную.
#   @END                                                          <- new
#   (blank separator paragraph)                                   <- new
#   @BEGIN BibliographySourceNotReferenced                        <- new
#   Указываемый источник не был использован в тексте.             <- new
#   @END                                                          <- original, pushed down

function Insert-ParagraphBeforeLastWithText($text) {
    # Inserts a brand new paragraph immediately before the current last
    # paragraph in the document, reusing that last paragraph's formatting,
    # and (optionally) fills it with $text.
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $insertPos = $lastPara.Range.Start
    $r = $d.Range($insertPos, $insertPos)
    $r.InsertParagraphBefore()

    $newParaIndex = $d.Paragraphs.Count - 1
    $newPara = $d.Paragraphs.Item($newParaIndex)
    if ($text) {
        $newPara.Range.InsertAfter($text)
    }
}

Insert-ParagraphBeforeLastWithText("@END")
Insert-ParagraphBeforeLastWithText("")
Insert-ParagraphBeforeLastWithText("@BEGIN BibliographySourceNotReferenced")
Insert-ParagraphBeforeLastWithText("Указываемый источник не был использован в тексте.")
